# The commit swaps the two embedded DrawingML themes:
#   ppt/theme/theme1.xml  (the slide master's theme, currently "Integral")
#     -> becomes the "Office Theme" color scheme
#   ppt/theme/theme2.xml  (the notes master's theme, currently "Office Theme")
#     -> becomes the "Integral" color scheme
#
# The font scheme (majorFont/minorFont) and format scheme (fills/lines/
# effects) are already byte-identical between the two themes, so only the
# twelve theme colours actually change. We drive the swap the same way a
# real PowerPoint COM automation script would: through
# ThemeColorScheme.Colors(i).RGB (a.k.a. .Item(i).RGB) on the master's Theme.

function Convert-HexToComRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    # PowerPoint's COM RGB values are packed 0x00BBGGRR, same as VBA's RGB().
    return $b * 65536 + $g * 256 + $r
}

$p = $ppt.ActivePresentation

# Target "Office Theme" colours (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# in ThemeColorScheme.Colors(1..12) order.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$masterColorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $masterColorScheme.Colors($i).RGB = Convert-HexToComRGB $officeThemeColors[$i - 1]
}

# Target "Integral" colours, applied the same way via the Notes Master's
# theme (ppt/theme/theme2.xml) so that part tracks the swap too.
$integralThemeColors = @(
    "000000",
    "FFFFFF",
    "455F51",
    "E3DED1",
    "99CB38",
    "63A537",
    "E6D024",
    "CC9700",
    "4EB3CF",
    "378DA6",
    "6B9F25",
    "B26B02"
)

$notesColorScheme = $p.NotesMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $integralThemeColors.Count; $i++) {
    $notesColorScheme.Colors($i).RGB = Convert-HexToComRGB $integralThemeColors[$i - 1]
}
